$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.481.23"
$ws.Range("E2").Value = "  -2.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.415.82"
$ws.Range("E3").Value = "  -2.23%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.10"
$ws.Range("E5").Value = "  -1.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.47"
$ws.Range("E6").Value = "  -3.73%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.156"
$ws.Range("E9").Value = "  +2.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  -1.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.326"
$ws.Range("E11").Value = "  -2.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.75"
$ws.Range("E12").Value = "  -1.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "67.419.31"
$ws.Range("E13").Value = "  -2.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000168"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.79"
$ws.Range("E15").Value = "  -3.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.31"
$ws.Range("E16").Value = "  -4.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "327.75"
$ws.Range("E17").Value = "  -4.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.78"
$ws.Range("E18").Value = "  -4.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.75"
$ws.Range("E19").Value = "  -1.56%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("E21").Value = "  -3.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.74"
$ws.Range("E22").Value = "  -2.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.59"
$ws.Range("E23").Value = "  -3.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.97"
$ws.Range("E24").Value = "  -2.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0₃0796"
$ws.Range("E25").Value = "  -3.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.97"
$ws.Range("E26").Value = "  -3.28%  "

$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "413.26"
$ws.Range("E28").Value = "  -6.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.11"
$ws.Range("E29").Value = "  -3.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.58"
$ws.Range("E30").Value = "  -2.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.28"
$ws.Range("E31").Value = "  +1.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.93"
$ws.Range("E32").Value = "  -0.74%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.68"
$ws.Range("E34").Value = "  -1.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.103"
$ws.Range("E35").Value = "  -5.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.291"
$ws.Range("E36").Value = "  -4.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.19"
$ws.Range("E37").Value = "  -6.30%  "

$ws.Range("E38").Value = "  -2.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.06"
$ws.Range("E39").Value = "  -4.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.96"
$ws.Range("E40").Value = "  -6.12%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "128.85"
$ws.Range("E41").Value = "  -3.60%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.28"
$ws.Range("E42").Value = "  -2.55%  "

$ws.Range("E43").Value = "  -1.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.472"
$ws.Range("E44").Value = "  -2.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.550"
$ws.Range("E45").Value = "  -2.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0909"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("E47").Value = "  -0.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.32"
$ws.Range("E48").Value = "  -8.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.39"
$ws.Range("E49").Value = "  -3.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0202"
$ws.Range("E50").Value = "  -2.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0425"
$ws.Range("E51").Value = "  -1.53%  "
